# TC10_BookedHotelPrice.xlsx edit
# - Adds a new "price" column (R) with value "AUD $ 250"
# - Updates the Check In / Check Out dates in row 2 (G2, H2) to the new dates
# - Updates the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header/value pair for the price column (column R)
$ws.Range("R1").Value = "price"
$ws.Range("R2").Value = "AUD $ 250"

# Update the booked check-in / check-out dates
$ws.Range("G2").Value = "29/05/2016"
$ws.Range("H2").Value = "30/05/2016"

# Update the selected cell to match the new state
$ws.Range("H8").Select()
